# Form the consolidated report: recompute the "Absent" column (H) from the
# "Real" column (E). A row is Absent (H=1) when its Real attendance count
# (E) is 0; otherwise the person attended and H=0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 3) { $lastRow = 21 }

for ($r = 3; $r -le $lastRow; $r++) {
    $realValue = $ws.Cells.Item($r, 5).Value2   # column E = Real
    if ($realValue -eq 0) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
